$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1.033209817776667
$ws.Cells.Item(2, 4).Value = 1.040364061206416
$ws.Cells.Item(2, 5).Value = 1.032405580217793
$ws.Cells.Item(2, 6).Value = 1.047653177649393
$ws.Cells.Item(2, 10).Value = 1.038335797327946
$ws.Cells.Item(2, 11).Value = 1.043146619765435
$ws.Cells.Item(2, 12).Value = 1.035210904742591
$ws.Cells.Item(2, 13).Value = 1.050415214237927
$ws.Cells.Item(2, 14).Value = 1.016343760030792

$ws.Cells.Item(3, 3).Value = 1.035205294500164
$ws.Cells.Item(3, 4).Value = 1.042231984745102
$ws.Cells.Item(3, 5).Value = 1.034140018643151
$ws.Cells.Item(3, 6).Value = 1.049703860060463
$ws.Cells.Item(3, 10).Value = 1.039968592958494
$ws.Cells.Item(3, 11).Value = 1.04482217148767
$ws.Cells.Item(3, 12).Value = 1.036751634166232
$ws.Cells.Item(3, 13).Value = 1.052274576819546
$ws.Cells.Item(3, 14).Value = 1.01692184552827

$ws.Cells.Item(4, 3).Value = 1.036492677772039
$ws.Cells.Item(4, 4).Value = 1.043437184137004
$ws.Cells.Item(4, 5).Value = 1.035258926093309
$ws.Cells.Item(4, 6).Value = 1.051027453646606
$ws.Cells.Item(4, 10).Value = 1.041021176213342
$ws.Cells.Item(4, 11).Value = 1.045902485317221
$ws.Cells.Item(4, 12).Value = 1.037744766336921
$ws.Cells.Item(4, 13).Value = 1.053474000180529
$ws.Cells.Item(4, 14).Value = 1.017293601880226

$ws.Cells.Item(5, 3).Value = 1.037033002415419
$ws.Cells.Item(5, 4).Value = 1.043943040265628
$ws.Cells.Item(5, 5).Value = 1.035728522832349
$ws.Cells.Item(5, 6).Value = 1.051583117563118
$ws.Cells.Item(5, 10).Value = 1.041462757539929
$ws.Cells.Item(5, 11).Value = 1.046355740361175
$ws.Cells.Item(5, 12).Value = 1.038161382251197
$ws.Cells.Item(5, 13).Value = 1.053977371825858
$ws.Cells.Item(5, 14).Value = 1.017449343421188

$ws.Cells.Item(6, 3).Value = 1.037123673732543
$ws.Cells.Item(6, 4).Value = 1.044027928877298
$ws.Cells.Item(6, 5).Value = 1.035807324319759
$ws.Cells.Item(6, 6).Value = 1.051676371226322
$ws.Cells.Item(6, 10).Value = 1.04153684733435
$ws.Cells.Item(6, 11).Value = 1.046431791159425
$ws.Cells.Item(6, 12).Value = 1.038231281835208
$ws.Cells.Item(6, 13).Value = 1.054061839995003
$ws.Cells.Item(6, 14).Value = 1.017475461354829

$ws.Cells.Item(7, 3).Value = 1.036499901085955
$ws.Cells.Item(7, 4).Value = 1.043443946563794
$ws.Cells.Item(7, 5).Value = 1.035265203953129
$ws.Cells.Item(7, 6).Value = 1.051034881467698
$ws.Cells.Item(7, 10).Value = 1.041027080251329
$ws.Cells.Item(7, 11).Value = 1.045908545277516
$ws.Cells.Item(7, 12).Value = 1.037750336678347
$ws.Cells.Item(7, 13).Value = 1.053480729630007
$ws.Cells.Item(7, 14).Value = 1.017295685038061

$ws.Cells.Item(8, 3).Value = 1.033885004760163
$ws.Cells.Item(8, 4).Value = 1.040996066642932
$ws.Cells.Item(8, 5).Value = 1.032992455666478
$ws.Cells.Item(8, 6).Value = 1.048346920193926
$ws.Cells.Item(8, 10).Value = 1.038888436951765
$ws.Cells.Item(8, 11).Value = 1.0437136958539
$ws.Cells.Item(8, 12).Value = 1.035732404095282
$ws.Cells.Item(8, 13).Value = 1.051044375663713
$ws.Cells.Item(8, 14).Value = 1.016539607745481

$ws.Cells.Item(9, 3).Value = 1.029246850724946
$ws.Cells.Item(9, 4).Value = 1.036655025604319
$ws.Cells.Item(9, 5).Value = 1.028960729416597
$ws.Cells.Item(9, 6).Value = 1.043583780703101
$ws.Cells.Item(9, 10).Value = 1.035088772547042
$ws.Cells.Item(9, 11).Value = 1.039815468083834
$ws.Cells.Item(9, 12).Value = 1.032146442497027
$ws.Cells.Item(9, 13).Value = 1.046721816072016
$ws.Cells.Item(9, 14).Value = 1.015189358393274

$ws.Cells.Item(10, 3).Value = 1.026132796579066
$ws.Cells.Item(10, 4).Value = 1.033741098287957
$ws.Cells.Item(10, 5).Value = 1.026253606906327
$ws.Cells.Item(10, 6).Value = 1.040388965633011
$ws.Cells.Item(10, 10).Value = 1.032533499730399
$ws.Cells.Item(10, 11).Value = 1.037194787113906
$ws.Cells.Item(10, 12).Value = 1.02973439348045
$ws.Cells.Item(10, 13).Value = 1.043818955894508
$ws.Cells.Item(10, 14).Value = 1.0142767017768

$ws.Cells.Item(11, 3).Value = 1.02477882706353
$ws.Cells.Item(11, 4).Value = 1.03247430891728
$ws.Cells.Item(11, 5).Value = 1.025076533615259
$ws.Cells.Item(11, 6).Value = 1.039000640693448
$ws.Cells.Item(11, 10).Value = 1.031421500692189
$ws.Cells.Item(11, 11).Value = 1.036054532196017
$ws.Cells.Item(11, 12).Value = 1.028684610337649
$ws.Cells.Item(11, 13).Value = 1.042556655902952
$ws.Cells.Item(11, 14).Value = 1.013878453048855

$ws.Cells.Item(12, 3).Value = 1.024275036930368
$ws.Cells.Item(12, 4).Value = 1.032002983224955
$ws.Cells.Item(12, 5).Value = 1.024638560003959
$ws.Cells.Item(12, 6).Value = 1.038484182377525
$ws.Cells.Item(12, 10).Value = 1.031007596704836
$ws.Cells.Item(12, 11).Value = 1.035630142264607
$ws.Cells.Item(12, 12).Value = 1.028293847671017
$ws.Cells.Item(12, 13).Value = 1.042086951809284
$ws.Cells.Item(12, 14).Value = 1.013730057304735

$ws.Cells.Item(13, 3).Value = 1.024383141274502
$ws.Cells.Item(13, 4).Value = 1.032104120081373
$ws.Cells.Item(13, 5).Value = 1.024732541443185
$ws.Cells.Item(13, 6).Value = 1.038594999873626
$ws.Cells.Item(13, 10).Value = 1.031096419737625
$ws.Cells.Item(13, 11).Value = 1.035721214142585
$ws.Cells.Item(13, 12).Value = 1.028377705356542
$ws.Cells.Item(13, 13).Value = 1.042187742914526
$ws.Cells.Item(13, 14).Value = 1.013761910028577

$ws.Cells.Item(14, 3).Value = 1.024737201444768
$ws.Cells.Item(14, 4).Value = 1.032435365131201
$ws.Cells.Item(14, 5).Value = 1.025040346147724
$ws.Cells.Item(14, 6).Value = 1.038957966018965
$ws.Cells.Item(14, 10).Value = 1.031387304905239
$ws.Cells.Item(14, 11).Value = 1.036019469445283
$ws.Cells.Item(14, 12).Value = 1.028652326765828
$ws.Cells.Item(14, 13).Value = 1.042517847098967
$ws.Cells.Item(14, 14).Value = 1.01386619621502

$ws.Cells.Item(15, 3).Value = 1.024955234008423
$ws.Cells.Item(15, 4).Value = 1.032639351487303
$ws.Cells.Item(15, 5).Value = 1.025229893850176
$ws.Cells.Item(15, 6).Value = 1.039181498183368
$ws.Cells.Item(15, 10).Value = 1.031566414425137
$ws.Cells.Item(15, 11).Value = 1.036203121201447
$ws.Cells.Item(15, 12).Value = 1.028821419880528
$ws.Cells.Item(15, 13).Value = 1.042721124426383
$ws.Cells.Item(15, 14).Value = 1.013930388052358

$ws.Cells.Item(16, 3).Value = 1.026222535333185
$ws.Cells.Item(16, 4).Value = 1.03382506248829
$ws.Cells.Item(16, 5).Value = 1.026331620682544
$ws.Cells.Item(16, 6).Value = 1.04048099744455
$ws.Cells.Item(16, 10).Value = 1.032607180461784
$ws.Cells.Item(16, 11).Value = 1.037270344456555
$ws.Cells.Item(16, 12).Value = 1.029803949482803
$ws.Cells.Item(16, 13).Value = 1.043902615730689
$ws.Cells.Item(16, 14).Value = 1.014303067002877

$ws.Cells.Item(17, 3).Value = 1.027015969617862
$ws.Cells.Item(17, 4).Value = 1.034567459787456
$ws.Cells.Item(17, 5).Value = 1.027021383557602
$ws.Cells.Item(17, 6).Value = 1.041294793619089
$ws.Cells.Item(17, 10).Value = 1.033258522782284
$ws.Cells.Item(17, 11).Value = 1.037938300148852
$ws.Cells.Item(17, 12).Value = 1.030418816083113
$ws.Cells.Item(17, 13).Value = 1.044642284737378
$ws.Cells.Item(17, 14).Value = 1.014536012922803

$ws.Cells.Item(18, 3).Value = 1.027478231641089
$ws.Cells.Item(18, 4).Value = 1.03500000313566
$ws.Cells.Item(18, 5).Value = 1.027423242214803
$ws.Cells.Item(18, 6).Value = 1.041768991645074
$ws.Cells.Item(18, 10).Value = 1.033637905997685
$ws.Cells.Item(18, 11).Value = 1.03832738008622
$ws.Cells.Item(18, 12).Value = 1.030776942719598
$ws.Cells.Item(18, 13).Value = 1.045073207538318
$ws.Cells.Item(18, 14).Value = 1.01467159143078

$ws.Cells.Item(19, 3).Value = 1.027635761231493
$ws.Cells.Item(19, 4).Value = 1.035147407869116
$ws.Cells.Item(19, 5).Value = 1.027560186920779
$ws.Cells.Item(19, 6).Value = 1.041930601136433
$ws.Cells.Item(19, 10).Value = 1.033767176070035
$ws.Cells.Item(19, 11).Value = 1.038459957592012
$ws.Cells.Item(19, 12).Value = 1.030898968036121
$ws.Cells.Item(19, 13).Value = 1.045220054771651
$ws.Cells.Item(19, 14).Value = 1.014717770427696

$ws.Cells.Item(20, 3).Value = 1.026930897155921
$ws.Cells.Item(20, 4).Value = 1.034487857882049
$ws.Cells.Item(20, 5).Value = 1.026947427148485
$ws.Cells.Item(20, 6).Value = 1.04120753038543
$ws.Cells.Item(20, 10).Value = 1.033188695298134
$ws.Cells.Item(20, 11).Value = 1.037866689545704
$ws.Cells.Item(20, 12).Value = 1.030352900111843
$ws.Cells.Item(20, 13).Value = 1.044562978573577
$ws.Cells.Item(20, 14).Value = 1.014511050593771

$ws.Cells.Item(21, 3).Value = 1.02463296368405
$ws.Cells.Item(21, 4).Value = 1.032337843569791
$ws.Cells.Item(21, 5).Value = 1.024949726405381
$ws.Cells.Item(21, 6).Value = 1.038851103099653
$ws.Cells.Item(21, 10).Value = 1.031301670356717
$ws.Cells.Item(21, 11).Value = 1.035931664315433
$ws.Cells.Item(21, 12).Value = 1.028571480603897
$ws.Cells.Item(21, 13).Value = 1.042420662766409
$ws.Cells.Item(21, 14).Value = 1.013835499531121

$ws.Cells.Item(22, 3).Value = 1.023183138755015
$ws.Cells.Item(22, 4).Value = 1.030981496454734
$ws.Cells.Item(22, 5).Value = 1.023689305324465
$ws.Cells.Item(22, 6).Value = 1.037365038888683
$ws.Cells.Item(22, 10).Value = 1.030110246246446
$ws.Cells.Item(22, 11).Value = 1.034710115551423
$ws.Cells.Item(22, 12).Value = 1.027446638298704
$ws.Cells.Item(22, 13).Value = 1.04106888987864
$ws.Cells.Item(22, 14).Value = 1.013408039319671

$ws.Cells.Item(23, 3).Value = 1.023952204531718
$ws.Cells.Item(23, 4).Value = 1.031700961772327
$ws.Cells.Item(23, 5).Value = 1.024357902465237
$ws.Cells.Item(23, 6).Value = 1.038153264580367
$ws.Cells.Item(23, 10).Value = 1.030742322659691
$ws.Cells.Item(23, 11).Value = 1.035358156535574
$ws.Cells.Item(23, 12).Value = 1.028043400484504
$ws.Cells.Item(23, 13).Value = 1.041785955696241
$ws.Cells.Item(23, 14).Value = 1.013634904183486

$ws.Cells.Item(24, 3).Value = 1.026969339391032
$ws.Cells.Item(24, 4).Value = 1.034523828052758
$ws.Cells.Item(24, 5).Value = 1.026980846312839
$ws.Cells.Item(24, 6).Value = 1.041246962352937
$ws.Cells.Item(24, 10).Value = 1.033220248981014
$ws.Cells.Item(24, 11).Value = 1.037899048923448
$ws.Cells.Item(24, 12).Value = 1.030382686291497
$ws.Cells.Item(24, 13).Value = 1.044598815201519
$ws.Cells.Item(24, 14).Value = 1.014522330906808

$ws.Cells.Item(25, 3).Value = 1.030449680348662
$ws.Cells.Item(25, 4).Value = 1.037780694888123
$ws.Cells.Item(25, 5).Value = 1.030006337236424
$ws.Cells.Item(25, 6).Value = 1.044818472727445
$ws.Cells.Item(25, 10).Value = 1.036074890452035
$ws.Cells.Item(25, 11).Value = 1.040827012393664
$ws.Cells.Item(25, 12).Value = 1.033077186371113
$ws.Cells.Item(25, 13).Value = 1.047842926655275
$ws.Cells.Item(25, 14).Value = 1.015540600153759
